$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.3270505502951835
$ws.Range("C2").Value = 0.05059594725676675
$ws.Range("D2").Value = 0.03169699048302732
$ws.Range("F2").Value = 0.6681628566188209
$ws.Range("G2").Value = 0.5082445335594912
$ws.Range("H2").Value = 0.6539383294400807
$ws.Range("I2").Value = 0.4951984570817061
$ws.Range("K2").Value = 0.3441850223563563
$ws.Range("N2").Value = 1.2816510421465
# Row 3
$ws.Range("B3").Value = 0.2900796426515626
$ws.Range("C3").Value = 0.04414566512869555
$ws.Range("D3").Value = 0.02928759997225683
$ws.Range("F3").Value = 0.6687410565119478
$ws.Range("G3").Value = 0.5103910393110809
$ws.Range("H3").Value = 0.6589369627200412
$ws.Range("I3").Value = 0.5004632671834415
$ws.Range("K3").Value = 0.3028698510027539
$ws.Range("N3").Value = 1.296030060003176
# Row 4
$ws.Range("B4").Value = 0.267393040291978
$ws.Range("C4").Value = 0.04017159623822408
$ws.Range("D4").Value = 0.02779628046827298
$ws.Range("F4").Value = 0.6695553835644148
$ws.Range("G4").Value = 0.5121063113227748
$ws.Range("H4").Value = 0.6623251112328674
$ws.Range("I4").Value = 0.503993618957022
$ws.Range("K4").Value = 0.2774933562772333
$ws.Range("N4").Value = 1.30534558168949
# Row 5
$ws.Range("B5").Value = 0.2581520103259152
$ws.Range("C5").Value = 0.03854875941664204
$ws.Range("D5").Value = 0.02718558653982228
$ws.Range("F5").Value = 0.6700026743331691
$ws.Range("G5").Value = 0.5129050647813003
$ws.Range("H5").Value = 0.6637860344121265
$ws.Range("I5").Value = 0.5055070875369267
$ws.Range("K5").Value = 0.2671505079258054
$ws.Range("N5").Value = 1.309264127776732
# Row 6
$ws.Range("B6").Value = 0.2566177974229618
$ws.Range("C6").Value = 0.03827908617215314
$ws.Range("D6").Value = 0.02708400295964708
$ws.Range("F6").Value = 0.6700839177704978
$ws.Range("G6").Value = 0.5130437187025549
$ws.Range("H6").Value = 0.6640334655372229
$ws.Range("I6").Value = 0.505762915509905
$ws.Range("K6").Value = 0.2654329964992428
$ws.Range("N6").Value = 1.309922191762039
# Row 7
$ws.Range("B7").Value = 0.267268395891108
$ws.Range("C7").Value = 0.04014972369053282
$ws.Range("D7").Value = 0.02778805641222704
$ws.Range("F7").Value = 0.6695609485166187
$ws.Range("G7").Value = 0.5121166798278765
$ws.Range("H7").Value = 0.662344488901411
$ws.Range("I7").Value = 0.5040137272362202
$ws.Range("K7").Value = 0.2773538753179139
$ws.Range("N7").Value = 1.305397932999004
# Row 8
$ws.Range("B8").Value = 0.3143004319447584
$ws.Range("C8").Value = 0.04837472890172023
$ws.Range("D8").Value = 0.03086873150284219
$ws.Range("F8").Value = 0.6682668573550714
$ws.Range("G8").Value = 0.5089021018886086
$ws.Range("H8").Value = 0.6555956671639294
$ws.Range("I8").Value = 0.4969519402101703
$ws.Range("K8").Value = 0.3299416649520879
$ws.Range("N8").Value = 1.286507903637403
# Row 9
$ws.Range("B9").Value = 0.4066210584593364
$ws.Range("C9").Value = 0.06439523828022686
$ws.Range("D9").Value = 0.03681391102919207
$ws.Range("F9").Value = 0.6693768658479939
$ws.Range("G9").Value = 0.5057576858300052
$ws.Range("H9").Value = 0.6448914853663439
$ws.Range("I9").Value = 0.4854683594988067
$ws.Range("K9").Value = 0.4329794520630514
$ws.Range("N9").Value = 1.253326318018303
# Row 10
$ws.Range("B10").Value = 0.4744873930038693
$ws.Range("C10").Value = 0.07609896872961031
$ws.Range("D10").Value = 0.04112201369503055
$ws.Range("F10").Value = 0.6724221181123724
$ws.Range("G10").Value = 0.5053838687625358
$ws.Range("H10").Value = 0.6385690382332712
$ws.Range("I10").Value = 0.4784759091175097
$ws.Range("K10").Value = 0.5086132295336085
$ws.Range("N10").Value = 1.23130077952462
# Row 11
$ws.Range("B11").Value = 0.5053666456857684
$ws.Range("C11").Value = 0.08140885855658553
$ws.Range("D11").Value = 0.04306864096704999
$ws.Range("F11").Value = 0.674293015719087
$ws.Range("G11").Value = 0.5056365480019309
$ws.Range("H11").Value = 0.6360275126153567
$ws.Range("I11").Value = 0.4756091944844805
$ws.Range("K11").Value = 0.5430034088892342
$ws.Range("N11").Value = 1.221791477238373
# Row 12
$ws.Range("B12").Value = 0.5170603260757503
$ws.Range("C12").Value = 0.08341750802300396
$ws.Range("D12").Value = 0.04380385689733401
$ws.Range("F12").Value = 0.6750713872840706
$ws.Range("G12").Value = 0.5057931875841319
$ws.Range("H12").Value = 0.6351132088956888
$ws.Range("I12").Value = 0.4745688810542461
$ws.Range("K12").Value = 0.5560233652712725
$ws.Range("N12").Value = 1.218263934464524
# Row 13
$ws.Range("B13").Value = 0.514541874897759
$ws.Range("C13").Value = 0.08298500318450408
$ws.Range("D13").Value = 0.04364560142440865
$ws.Range("F13").Value = 0.6749006408939024
$ws.Range("G13").Value = 0.505756738785621
$ws.Range("H13").Value = 0.6353079804906088
$ws.Range("I13").Value = 0.4747909175890399
$ws.Range("K13").Value = 0.5532194206953136
$ws.Range("N13").Value = 1.219020387498244
# Row 14
$ws.Range("B14").Value = 0.5063286874489279
$ws.Range("C14").Value = 0.08157415353811359
$ws.Range("D14").Value = 0.04312916651618082
$ws.Range("F14").Value = 0.6743556514258557
$ws.Range("G14").Value = 0.5056482122423915
$ws.Range("H14").Value = 0.6359513279431468
$ws.Range("I14").Value = 0.4755227000133715
$ws.Range("K14").Value = 0.5440746290095717
$ws.Range("N14").Value = 1.221499792335724
# Row 15
$ws.Range("B15").Value = 0.5012979055433675
$ws.Range("C15").Value = 0.08070969367290104
$ws.Range("D15").Value = 0.04281258279205957
$ws.Range("F15").Value = 0.6740309353113361
$ws.Range("G15").Value = 0.5055896797457109
$ws.Range("H15").Value = 0.636351663181415
$ws.Range("I15").Value = 0.4759768327203737
$ws.Range("K15").Value = 0.5384727913018139
$ws.Range("N15").Value = 1.223028062951936
# Row 16
$ws.Range("B16").Value = 0.4724694445309012
$ws.Range("C16").Value = 0.07575166508121356
$ws.Range("D16").Value = 0.040994529110435
$ws.Range("F16").Value = 0.6723096284446228
$ws.Range("G16").Value = 0.5053758763248339
$ws.Range("H16").Value = 0.6387418659677309
$ws.Range("I16").Value = 0.4786695838426382
$ws.Range("K16").Value = 0.5063653807987407
$ws.Range("N16").Value = 1.231932506378683
# Row 17
$ws.Range("B17").Value = 0.4547853885378572
$ws.Range("C17").Value = 0.07270640425633701
$ws.Range("D17").Value = 0.03987581651783501
$ws.Range("F17").Value = 0.6713780892169581
$ws.Range("G17").Value = 0.5053531056545921
$ws.Range("H17").Value = 0.640293876322346
$ws.Range("I17").Value = 0.4804020196211738
$ws.Range("K17").Value = 0.4866640292147224
$ws.Range("N17").Value = 1.237525812517639
# Row 18
$ws.Range("B18").Value = 0.4446146341076656
$ws.Range("C18").Value = 0.07095351754651347
$ws.Range("D18").Value = 0.03923112716699251
$ws.Range("F18").Value = 0.6708879947775372
$ws.Range("G18").Value = 0.5053797880267581
$ws.Range("H18").Value = 0.6412180443395812
$ws.Range("I18").Value = 0.4814280393927675
$ws.Range("K18").Value = 0.4753308596768591
$ws.Range("N18").Value = 1.240790955743599
# Row 19
$ws.Range("B19").Value = 0.4411711195556904
$ws.Range("C19").Value = 0.07035979274053261
$ws.Range("D19").Value = 0.03901263541438027
$ws.Range("F19").Value = 0.6707299043301802
$ws.Range("G19").Value = 0.5053956492654805
$ws.Range("H19").Value = 0.6415363605471356
$ws.Range("I19").Value = 0.4817805083591225
$ws.Range("K19").Value = 0.471493407130481
$ws.Range("N19").Value = 1.241904723275205
# Row 20
$ws.Range("B20").Value = 0.4566678235927668
$ws.Range("C20").Value = 0.07303071581017662
$ws.Range("D20").Value = 0.03999503352826395
$ws.Range("F20").Value = 0.6714725225958915
$ws.Range("G20").Value = 0.5053514113912598
$ws.Range("H20").Value = 0.6401254027669125
$ws.Range("I20").Value = 0.4802145381376555
$ws.Range("K20").Value = 0.4887614282232278
$ws.Range("N20").Value = 1.236925425105252
# Row 21
$ws.Range("B21").Value = 0.5087410939075028
$ws.Range("C21").Value = 0.08198861140895985
$ws.Range("D21").Value = 0.04328090864311207
$ws.Range("F21").Value = 0.6745138302976699
$ws.Range("G21").Value = 0.5056784334775841
$ws.Range("H21").Value = 0.6357610552933437
$ws.Range("I21").Value = 0.4753065291729968
$ws.Range("K21").Value = 0.5467607580099241
$ws.Range("N21").Value = 1.220769537896242
# Row 22
$ws.Range("B22").Value = 0.5427759271308616
$ws.Range("C22").Value = 0.08783090561121298
$ws.Range("D22").Value = 0.04541715307092886
$ws.Range("F22").Value = 0.6769089842528757
$ws.Range("G22").Value = 0.5062475331356069
$ws.Range("H22").Value = 0.6331891656805766
$ws.Range("I22").Value = 0.4723626406440395
$ws.Range("K22").Value = 0.584649754433002
$ws.Range("N22").Value = 1.210638727454022
# Row 23
$ws.Range("B23").Value = 0.5246109034906112
$ws.Range("C23").Value = 0.08471389662042839
$ws.Range("D23").Value = 0.04427804270173397
$ws.Range("F23").Value = 0.6755933362718949
$ws.Range("G23").Value = 0.5059112223916316
$ws.Range("H23").Value = 0.634536167879773
$ws.Range("I23").Value = 0.4739096906049554
$ws.Range("K23").Value = 0.5644294160951802
$ws.Range("N23").Value = 1.21600655158818
# Row 24
$ws.Range("B24").Value = 0.4558167869823819
$ws.Range("C24").Value = 0.07288410117556055
$ws.Range("D24").Value = 0.03994114026987461
$ws.Range("F24").Value = 0.6714296876671639
$ws.Range("G24").Value = 0.5053520534886928
$ws.Range("H24").Value = 0.6402014703050867
$ws.Range("I24").Value = 0.4802992050131252
$ws.Range("K24").Value = 0.487813214648753
$ws.Range("N24").Value = 1.237196706139013
# Row 25
$ws.Range("B25").Value = 0.3816379106933994
$ws.Range("C25").Value = 0.06007293070155129
$ws.Range("D25").Value = 0.03521599465105396
$ws.Range("F25").Value = 0.6686854440352263
$ws.Range("G25").Value = 0.506269028963807
$ws.Range("H25").Value = 0.6475164323815932
$ws.Range("I25").Value = 0.4883215157260317
$ws.Range("K25").Value = 0.4051158555210748
$ws.Range("N25").Value = 1.261889531070981

Write-Host "Applied 216 value updates"